$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1177
$ws.Range("J17").Value = 1177
$ws.Range("L17").Value = 3531
$ws.Range("N17").Value = -3867
# Row 64
$ws.Range("H64").Value = 3944.9
$ws.Range("I64").Value = 3938.7778
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 3938.7778
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -3690.7778
$ws.Range("N64").Value = -4496
# Row 67
$ws.Range("H67").Value = 3944.9
$ws.Range("I67").Value = 3938.7778
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 3938.7778
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -3080.7778
$ws.Range("N67").Value = -5716
# Row 76
$ws.Range("H76").Value = 3110.606
$ws.Range("I76").Value = 3044
$ws.Range("K76").Value = 3044
$ws.Range("M76").Value = -2729
# Row 79
$ws.Range("H79").Value = 3110.606
$ws.Range("I79").Value = 3044
$ws.Range("K79").Value = 3044
$ws.Range("M79").Value = -1952
# Row 112
$ws.Range("H112").Value = 1517.2222
$ws.Range("J112").Value = 1517.2222
$ws.Range("L112").Value = 4551.6666
$ws.Range("N112").Value = -6767.6666
# Row 116
$ws.Range("H116").Value = 2418.5454
$ws.Range("I116").Value = 2151.25
$ws.Range("J116").Value = 2571.2856
$ws.Range("K116").Value = 2151.25
$ws.Range("L116").Value = 2571.2856
$ws.Range("M116").Value = 1290.75
$ws.Range("N116").Value = -9455.285599999999
# Row 125
$ws.Range("H125").Value = 1860.238
$ws.Range("I125").Value = 795.8
$ws.Range("J125").Value = 2192.875
$ws.Range("K125").Value = 7162.2
$ws.Range("L125").Value = 19735.875
$ws.Range("M125").Value = -4702.2
$ws.Range("N125").Value = -24655.875
# Row 129
$ws.Range("H129").Value = 531.75
$ws.Range("J129").Value = 925
$ws.Range("L129").Value = 2775
$ws.Range("N129").Value = -12775
# Row 138
$ws.Range("H138").Value = 2101.878
$ws.Range("I138").Value = 1841.4474
$ws.Range("J138").Value = 5400.6665
$ws.Range("K138").Value = 5524.3422
$ws.Range("L138").Value = 16201.9995
$ws.Range("M138").Value = -384.3422
$ws.Range("N138").Value = -26481.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1535.5625
$ws.Range("I105").Value = 1529.9231
$ws.Range("J105").Value = 1560
$ws.Range("K105").Value = 1529.9231
$ws.Range("L105").Value = 1560
$ws.Range("M105").Value = 217.0769
$ws.Range("N105").Value = -5054
# Row 134
$ws.Range("H134").Value = 34627
$ws.Range("I134").Value = 43557.19
$ws.Range("J134").Value = 1457.7142
$ws.Range("K134").Value = 130671.57
$ws.Range("L134").Value = 4373.142599999999
$ws.Range("M134").Value = -128136.57
$ws.Range("N134").Value = -9443.142599999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1588.0927
$ws.Range("I31").Value = 1298.871
$ws.Range("J31").Value = 1977.9131
$ws.Range("K31").Value = 1298.871
$ws.Range("L31").Value = 1977.9131
$ws.Range("M31").Value = -1003.871
$ws.Range("N31").Value = -2567.9131
# Row 34
$ws.Range("H34").Value = 1588.0927
$ws.Range("I34").Value = 1298.871
$ws.Range("J34").Value = 1977.9131
$ws.Range("K34").Value = 1298.871
$ws.Range("L34").Value = 1977.9131
$ws.Range("M34").Value = -1096.871
$ws.Range("N34").Value = -2381.9131
# Row 51
$ws.Range("H51").Value = 10832.538
$ws.Range("J51").Value = 11644.417
$ws.Range("L51").Value = 11644.417
$ws.Range("N51").Value = -13116.417
# Row 58
$ws.Range("H58").Value = 2451.506
$ws.Range("I58").Value = 1058.825
$ws.Range("J58").Value = 3747.0232
$ws.Range("K58").Value = 1058.825
$ws.Range("L58").Value = 3747.0232
$ws.Range("M58").Value = -855.825
$ws.Range("N58").Value = -4153.0232
# Row 61
$ws.Range("H61").Value = 10832.538
$ws.Range("J61").Value = 11644.417
$ws.Range("L61").Value = 11644.417
$ws.Range("N61").Value = -12340.417
# Row 62
$ws.Range("H62").Value = 2779.375
$ws.Range("I62").Value = 2353.3333
$ws.Range("J62").Value = 3327.1428
$ws.Range("K62").Value = 2353.3333
$ws.Range("L62").Value = 3327.1428
$ws.Range("M62").Value = -1729.3333
$ws.Range("N62").Value = -4575.1428
# Row 65
$ws.Range("H65").Value = 2779.375
$ws.Range("I65").Value = 2353.3333
$ws.Range("J65").Value = 3327.1428
$ws.Range("K65").Value = 11766.6665
$ws.Range("L65").Value = 16635.714
$ws.Range("M65").Value = -8646.666499999999
$ws.Range("N65").Value = -22875.714
# Row 136
$ws.Range("H136").Value = 2451.506
$ws.Range("I136").Value = 1058.825
$ws.Range("J136").Value = 3747.0232
$ws.Range("K136").Value = 3176.475
$ws.Range("L136").Value = 11241.0696
$ws.Range("M136").Value = -626.4750000000004
$ws.Range("N136").Value = -16341.0696

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 74
$ws.Range("H74").Value = 8416.666999999999
$ws.Range("I74").Value = 500
$ws.Range("J74").Value = 10000
$ws.Range("K74").Value = 1500
$ws.Range("L74").Value = 30000
$ws.Range("M74").Value = -439
$ws.Range("N74").Value = -32122
# Row 77
$ws.Range("H77").Value = 8416.666999999999
$ws.Range("I77").Value = 500
$ws.Range("J77").Value = 10000
$ws.Range("K77").Value = 4500
$ws.Range("L77").Value = 90000
$ws.Range("M77").Value = 804
$ws.Range("N77").Value = -100608
# Row 131
$ws.Range("H131").Value = 2648.2878
$ws.Range("I131").Value = 20243.334
$ws.Range("J131").Value = 1810.4286
$ws.Range("K131").Value = 60730.00199999999
$ws.Range("L131").Value = 5431.2858
$ws.Range("M131").Value = -55690.00199999999
$ws.Range("N131").Value = -15511.2858

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 13104.454
$ws.Range("I70").Value = 16143.625
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 16143.625
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -15873.625
$ws.Range("N70").Value = -5540
# Row 73
$ws.Range("H73").Value = 13104.454
$ws.Range("I73").Value = 16143.625
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 16143.625
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -15207.625
$ws.Range("N73").Value = -6872
# Row 80
$ws.Range("H80").Value = 12145.454
$ws.Range("J80").Value = 17107.143
$ws.Range("L80").Value = 17107.143
$ws.Range("N80").Value = -19103.143
# Row 83
$ws.Range("H83").Value = 12145.454
$ws.Range("J83").Value = 17107.143
$ws.Range("L83").Value = 85535.715
$ws.Range("N83").Value = -95519.715

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 87
$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52246
# Row 90
$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -161232
# Row 122
$ws.Range("H122").Value = 2510.889
$ws.Range("I122").Value = 1849.5
$ws.Range("J122").Value = 3040
$ws.Range("K122").Value = 5548.5
$ws.Range("L122").Value = 9120
$ws.Range("M122").Value = -3098.5
$ws.Range("N122").Value = -14020
# Row 136
$ws.Range("H136").Value = 2426.7368
$ws.Range("I136").Value = 1478.6666
$ws.Range("J136").Value = 3280
$ws.Range("K136").Value = 4435.9998
$ws.Range("L136").Value = 9840
$ws.Range("M136").Value = -1885.9998
$ws.Range("N136").Value = -14940

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2171.7144
$ws.Range("I81").Value = 1520.4
$ws.Range("J81").Value = 3800
$ws.Range("K81").Value = 3040.8
$ws.Range("L81").Value = 7600
$ws.Range("M81").Value = -1979.8
$ws.Range("N81").Value = -9722
# Row 84
$ws.Range("H84").Value = 2171.7144
$ws.Range("I84").Value = 1520.4
$ws.Range("J84").Value = 3800
$ws.Range("K84").Value = 15204
$ws.Range("L84").Value = 38000
$ws.Range("M84").Value = -9900
$ws.Range("N84").Value = -48608
# Row 126
$ws.Range("H126").Value = 2422.3103
$ws.Range("I126").Value = 2157.1667
$ws.Range("J126").Value = 3695
$ws.Range("K126").Value = 6471.500100000001
$ws.Range("L126").Value = 11085
$ws.Range("M126").Value = -4001.500100000001
$ws.Range("N126").Value = -16025
